$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column Q (years 2009..2021 already present through P) with the
# new 2022 figures, re-using the existing column-P cell formatting so no
# new style entries are created (matches P's per-row style exactly).
$ws.Range("P3:P25").Copy()
$ws.Range("Q3:Q25").PasteSpecial(-4122)

# Header: new year value
$ws.Range("Q4").Value = 2022

# Data values for 2022
$ws.Range("Q5").Value = 8725
$ws.Range("Q7").Value = 8347
$ws.Range("Q8").Value = 378

# Not-yet-available figures ("…") for the remaining breakdown rows
$ws.Range("Q10:Q25").Value = "…"

# Move the active selection to Q3, matching the saved view state
$ws.Range("Q3").Select()
